$d = $word.ActiveDocument

# Change 1: table row 1 (1-based), column 1 -> "USE CASE #01"
$cell1 = $d.Tables.Item(1).Cell(1, 1)
$r1 = $cell1.Range
$r1.Find.Execute(".", $true, $false, $false, $false, $false, $true, 0, $false, "USE CASE #01", 1)

# Change 2: table row 40 (1-based), column 3 -> merge runs into " Inserisce un nuovo indirizzo."
$cell2 = $d.Tables.Item(1).Cell(40, 3)
$r2 = $cell2.Range
$r2.Find.Execute(" Inserisce un nuovo indirizzo.", $true, $false, $false, $false, $false, $true, 0, $false, " Inserisce un nuovo indirizzo.", 1)

# Change 3: table row 26 (1-based), column 3 -> merge runs into "creazioneAnnuncioNext05."
$cell3 = $d.Tables.Item(1).Cell(26, 3)
$r3 = $cell3.Range
$r3.Find.Execute("creazioneAnnuncioNext05.", $true, $false, $false, $false, $false, $true, 0, $false, "creazioneAnnuncioNext05.", 1)
